$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has a single "StoresList" / "Philadelphia, PA, USA"
# column. The new layout adds a "url" column in front of it, turning the
# old single column into a "location" column. Inserting a new column A
# shifts the existing column (values, styles, and width) over to column B
# intact, which is exactly the "location" column the diff expects.
$ws.Columns.Item(1).Insert()

# A2 is a regular bordered data cell - copy the border-only style now
# carried by B2 (the old A2) so it reuses the existing style definition
# instead of minting a new one.
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Populate the new column A ("url" header + link). Column order (A then B)
# matches the original authoring order so new shared-string entries land
# in the expected slots.
$ws.Range("A1").Value = "url"
$ws.Range("A2").Value = "https://dev01.ip.wawa.com/commerce/ui/"

# Fix up the header row text: old column now reads "location" instead of
# "StoresList".
$ws.Range("B1").Value = "location"

# A1 ("url") is a bold header with no border - a style not previously used
# in this sheet, so Excel mints a new cellXfs entry for it.
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Borders.LineStyle = -4142

# Column A width (engine quantizes ColumnWidth to 1/6-character steps, so
# use the closest reachable input to the desired 39.66-character width).
$ws.Columns.Item(1).ColumnWidth = 38.8333

# Selection ends up parked on A10, as in the authored workbook.
$ws.Range("A10").Select()
